$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 'Excellent app'
$ws.Range("A6").Value = 'Great work'
$ws.Range("A7").Value = 'Good job'
$ws.Range("A8").Value = 'Good job'
$ws.Range("A9").Value = 'Nice work'
$ws.Range("A10").Value = 'Nicely done'
$ws.Range("A11").Value = 'Nicely done'
$ws.Range("A12").Value = 'Good work'
$ws.Range("A13").Value = 'Great work'
$ws.Range("A14").Value = 'Good work'
$ws.Range("A15").Value = 'Good work team'
$ws.Range("A16").Value = 'Good work team'
$ws.Range("A17").Value = 'GOod'
$ws.Range("A18").Value = 'Good'
$ws.Range("A19").Value = 'Great'
$ws.Range("A20").Value = 'Great'
$ws.Range("A21").Value = 'Great'
$ws.Range("A22").Value = 'Good work on backtesting'
$ws.Range("A23").Value = 'Good work'
$ws.Range("A24").Value = 'Peers food'
$ws.Range("A25").Value = 'Back'
$ws.Range("A26").Value = 'Good job'
$ws.Range("A27").Value = 'wow'
$ws.Range("A28").Value = 'nice'
$ws.Range("A29").Value = 'fod'
$ws.Range("A30").Value = 'ads'
$ws.Range("A31").Value = 'Good job'
$ws.Range("A32").Value = 'nicely done'
$ws.Range("A33").Value = 'wow'
$ws.Range("A34").Value = 'good'
$ws.Range("A35").Value = 'wow'
$ws.Range("A36").Value = 'technicals'
$ws.Range("A37").Value = 'good'
$ws.Range("A38").Value = 'asd'
$ws.Range("A39").Value = 'asd'
$ws.Range("A40").Value = 'asd'
$ws.Range("A41").Value = 'very good'
